# Fruta / hortaliza, semanal
# Update the "Fecha" (D) values: rows 2-4 move from 2021-04-09 (44295) to
# 2021-04-08 (44294), and rows 5-7 move from 2021-04-08 (44294) to
# 2021-04-09 (44295). Also swap the "Volumen" (M) values for rows 3 and 6
# so they stay consistent with their (now swapped) date/quality grouping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44294
$ws.Range("D3").Value = 44294
$ws.Range("D4").Value = 44294

$ws.Range("D5").Value = 44295
$ws.Range("D6").Value = 44295
$ws.Range("D7").Value = 44295

$ws.Range("M3").Value = 240
$ws.Range("M6").Value = 200
